# issue #5: stock data from json to db
# Adds a "category" column (after property_category) and "source_file" /
# "index" columns (at the end) to the 股票 (stock) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at I: shifts old I(date)/J(legislator_name)/K(legislator_id)
# one column to the right -> J/K/L, freeing up column I for "category".
$ws.Columns.Item(9).Insert()

# ---- Header row (row 1) ----
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Copy header style (bold, bordered, centered) from an existing header cell
# onto the newly added header cells.
$ws.Cells.Item(1, 8).Copy() | Out-Null
$ws.Cells.Item(1, 9).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 8).Copy() | Out-Null
$ws.Cells.Item(1, 13).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 8).Copy() | Out-Null
$ws.Cells.Item(1, 14).PasteSpecial(-4122) | Out-Null

# ---- Data rows (2-4) ----
$ws.Cells.Item(2, 9).Value = "normal"
$ws.Cells.Item(3, 9).Value = "normal"
$ws.Cells.Item(4, 9).Value = "normal"

$ws.Cells.Item(2, 13).Value = "tmp43441"
$ws.Cells.Item(3, 13).Value = "tmp43441"
$ws.Cells.Item(4, 13).Value = "tmp43441"

$ws.Cells.Item(2, 14).Value = 72
$ws.Cells.Item(3, 14).Value = 73
$ws.Cells.Item(4, 14).Value = 74

# Copy data-row style (bordered, non-bold) from an existing data cell onto
# the newly added data cells.
$ws.Cells.Item(2, 8).Copy() | Out-Null
$ws.Cells.Item(2, 9).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(3, 8).Copy() | Out-Null
$ws.Cells.Item(3, 9).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(4, 8).Copy() | Out-Null
$ws.Cells.Item(4, 9).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(2, 8).Copy() | Out-Null
$ws.Cells.Item(2, 13).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(3, 8).Copy() | Out-Null
$ws.Cells.Item(3, 13).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(4, 8).Copy() | Out-Null
$ws.Cells.Item(4, 13).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(2, 8).Copy() | Out-Null
$ws.Cells.Item(2, 14).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(3, 8).Copy() | Out-Null
$ws.Cells.Item(3, 14).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(4, 8).Copy() | Out-Null
$ws.Cells.Item(4, 14).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

Write-Output "done"
